# Diagrams.xlsx re-layout:
#  - insert a new "Description" column (new column B)
#  - insert a new top row carrying two group headers:
#      "For Physics Group"  (over Diagram/Description)
#      "For Graphics Group" (over Date/Creator/Status/File, merged + centered)
#  - give the (now) second row (the original bold header row) a light-gray
#    fill to match the new group-header banner above it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edits -----------------------------------------------------

# New column for "Description", inserted before the old Date column (old B).
# Everything from the old B:E shifts right to C:F.
$ws.Columns("B:B").Insert()

# New top row for the group headers. Old row 1 (headers) -> row 2,
# old row 2 (data) -> row 3.
$ws.Rows("1:1").Insert()

# --- content ---------------------------------------------------------------

# New "Description" header, alongside the existing header row (now row 2).
$ws.Range("B2").Value = "Description"

# New group-header banner row (row 1).
$ws.Range("A1").Value = "For Physics Group"
$ws.Range("C1").Value = "For Graphics Group"

# --- formatting --------------------------------------------------------------

# "For Physics Group" banner: yellow fill over A1:B1.
$ws.Range("A1:B1").Interior.Color = 65535

# "For Graphics Group" banner: orange fill + centered, merged over C1:F1.
$ws.Range("C1:F1").Interior.Color = 49407
$ws.Range("C1:F1").HorizontalAlignment = -4108
$ws.Range("C1:F1").Merge()

# Header row (now row 2) gets a light-gray fill behind the existing
# bold font + bottom border.
$ws.Range("A2:F2").Interior.Color = 14277081

# Restore the user's last selection.
$ws.Range("B10").Select() | Out-Null
